$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.973.72"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "2.968.41"
$ws.Range("E3").Value = "  +1.23%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "352.76"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.00"
$ws.Range("E6").Value = "  -3.75%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.558"
$ws.Range("E7").Value = "  -0.40%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.613"
$ws.Range("E9").Value = "  -1.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.29"
$ws.Range("E10").Value = "  -2.80%  "
$ws.Range("E11").Value = "  +1.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0852"
$ws.Range("E12").Value = "  -4.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.09"
$ws.Range("E13").Value = "  -4.72%  "
$ws.Range("D14").Value = "3.422.68"
$ws.Range("E14").Value = "  +0.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.60"
$ws.Range("E15").Value = "  -2.34%  "
$ws.Range("D16").Value = "2.961.88"
$ws.Range("E16").Value = "  +0.81%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.994"
$ws.Range("E17").Value = "  +0.90%  "
$ws.Range("D18").Value = "51.899.93"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.43"
$ws.Range("E19").Value = "  +3.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.45"
$ws.Range("E20").Value = "  -2.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.54"
$ws.Range("E21").Value = "  -4.83%  "
$ws.Range("D22").Value = "0.0₃0970"
$ws.Range("E22").Value = "  -1.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.40"
$ws.Range("E23").Value = "  -2.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "263.57"
$ws.Range("E24").Value = "  -1.87%  "
$ws.Range("E25").Value = "  -2.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.176"
$ws.Range("E26").Value = "  -2.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.72"
$ws.Range("E27").Value = "  -0.77%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.41"
$ws.Range("E29").Value = "  +1.95%  "
$ws.Range("E30").Value = "  +1.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "10.30"
$ws.Range("E31").Value = "  -2.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.06"
$ws.Range("E32").Value = "  -3.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "36.11"
$ws.Range("E33").Value = "  -2.84%  "
$ws.Range("E34").Value = "  -4.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.61"
$ws.Range("E35").Value = "  -4.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0435"
$ws.Range("E36").Value = "  -4.37%  "
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.18"
$ws.Range("E38").Value = "  -4.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.82"
$ws.Range("E39").Value = "  -4.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.96"
$ws.Range("E40").Value = "  -4.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.71"
$ws.Range("E41").Value = "  +0.76%  "
$ws.Range("E42").Value = "  -1.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "123.09"
$ws.Range("E43").Value = "  +10.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.45"
$ws.Range("E44").Value = "  -3.11%  "
$ws.Range("E45").Value = "  -3.81%  "
$ws.Range("D46").Value = "2.113.34"
$ws.Range("E46").Value = "  -2.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.34"
$ws.Range("E47").Value = "  -4.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.32"
$ws.Range("E48").Value = "  -8.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.237"
$ws.Range("E49").Value = "  -4.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0339"
$ws.Range("E50").Value = "  -3.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.924"
$ws.Range("E51").Value = "  -2.17%  "
